# Generate Report for Handback
# ------------------------------------------------------------------
# The localization run finished and the target files are back "in sync"
# with en-US, so:
#   1. Every "Ready for handoff" status cell (Overview summary columns
#      + the per-locale Status column) becomes "Handed back: in sync
#      with en-US".
#   2. Each locale sheet (zh-cn, de-de) gets its "Latest Target File"
#      and "Latest Handback File" columns (F/G) populated with
#      hyperlinks to the source markdown file and the handback .xlf
#      file, for both data rows.
#   3. The "Latest Handback DateTime" column (H) is stamped with the
#      real handback time for each locale.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# ---- Overview sheet: zh-cn / de-de summary status columns (B, C) ----
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("B2").Value = $statusText
$ovw.Range("C2").Value = $statusText
$ovw.Range("B3").Value = $statusText
$ovw.Range("C3").Value = $statusText

# ---- locale detail sheets ----
$mdFile      = "f64fd7e2-04ab-449d-ada8-fdaffe3197cb.md"
$mdUrl       = "https://github.com/OpenLocalizationTest/oltest/blob/5ebae9d5bdfd8444ef3bb8d6ed4fd3816c76eabb/e2e/f64fd7e2-04ab-449d-ada8-fdaffe3197cb.md"

$locales = @(
    @{ Sheet = "zh-cn"; XlfFile = "f64fd7e2-04ab-449d-ada8-fdaffe3197cb.c4ef44521985cc8052aa5530f95c3ba80f4971b4.zh-cn.xlf";
       XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1441cde99ce9f542daf8adebc9cae070a6617c16/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/f64fd7e2-04ab-449d-ada8-fdaffe3197cb.c4ef44521985cc8052aa5530f95c3ba80f4971b4.zh-cn.xlf";
       HandbackTime = "2016-03-13 21:13:51" },
    @{ Sheet = "de-de"; XlfFile = "f64fd7e2-04ab-449d-ada8-fdaffe3197cb.c4ef44521985cc8052aa5530f95c3ba80f4971b4.de-de.xlf";
       XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e8930a5d851a81fa0429020f22d5bb626aa55757/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/f64fd7e2-04ab-449d-ada8-fdaffe3197cb.c4ef44521985cc8052aa5530f95c3ba80f4971b4.de-de.xlf";
       HandbackTime = "2016-03-13 21:13:57" }
)

foreach ($loc in $locales) {
    $ws = $wb.Worksheets.Item($loc.Sheet)

    # Status column (C) -> handed back
    $ws.Range("C2").Value = $statusText
    $ws.Range("C3").Value = $statusText

    # Latest Target File (F) + Latest Handback File (G), rows 2 & 3
    $ws.Hyperlinks.Add($ws.Range("F2"), $mdUrl, "", "", $mdFile)
    $ws.Hyperlinks.Add($ws.Range("G2"), $loc.XlfUrl, "", "", $loc.XlfFile)
    $ws.Hyperlinks.Add($ws.Range("F3"), $mdUrl, "", "", $mdFile)
    $ws.Hyperlinks.Add($ws.Range("G3"), $loc.XlfUrl, "", "", $loc.XlfFile)

    # Apply the same visual style the other hyperlink columns use
    $ws.Range("F2").Style = "HyperLink"
    $ws.Range("G2").Style = "HyperLink"
    $ws.Range("F3").Style = "HyperLink"
    $ws.Range("G3").Style = "HyperLink"

    # Latest Handback DateTime (H)
    $ws.Range("H2").Value = $loc.HandbackTime
    $ws.Range("H3").Value = $loc.HandbackTime
}
